# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File", "Latest Handback
#    File" and "Latest Handback DateTime" columns populated, with a new
#    hyperlink added on the "Latest Target File" cell
#  - a handful of columns are widened so the new content is readable

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5523676afec87277645b858897512a4cb90f012e/e2e/"
$targetMd = "a4685b61-dad9-48df-8bab-27e8dccbf679.md"
$otherMd  = "ffff5f93a686-4938-46d8-a668-d4541e460e27.md"

# ---------------------------------------------------------------------
# 1. Overview sheet: mark both language rows as handed back
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item(1)
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item(2)

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $baseUrl + $targetMd, [Type]::Missing, [Type]::Missing, $targetMd)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $baseUrl + $targetMd, [Type]::Missing, [Type]::Missing, $targetMd)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $baseUrl + $otherMd, [Type]::Missing, [Type]::Missing, $otherMd)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $baseUrl + $targetMd, [Type]::Missing, [Type]::Missing, $targetMd)

$zhcn.Range("J2").Value = "a4685b61-dad9-48df-8bab-27e8dccbf679.54df1e724f8c95c1c69dbb290e408c3b949889c7.zh-cn.xlf"
$zhcn.Range("J3").Value = "a4685b61-dad9-48df-8bab-27e8dccbf679.54df1e724f8c95c1c69dbb290e408c3b949889c7.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-01 03:14:03"
$zhcn.Range("K3").Value = "2016-09-01 03:14:03"

$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item(3)

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $baseUrl + $targetMd, [Type]::Missing, [Type]::Missing, $targetMd)
$dede.Hyperlinks.Add($dede.Range("I2"), $baseUrl + $targetMd, [Type]::Missing, [Type]::Missing, $targetMd)
$dede.Hyperlinks.Add($dede.Range("A3"), $baseUrl + $otherMd, [Type]::Missing, [Type]::Missing, $otherMd)
$dede.Hyperlinks.Add($dede.Range("I3"), $baseUrl + $targetMd, [Type]::Missing, [Type]::Missing, $targetMd)

$dede.Range("J2").Value = "a4685b61-dad9-48df-8bab-27e8dccbf679.54df1e724f8c95c1c69dbb290e408c3b949889c7.de-de.xlf"
$dede.Range("J3").Value = "a4685b61-dad9-48df-8bab-27e8dccbf679.54df1e724f8c95c1c69dbb290e408c3b949889c7.de-de.xlf"
$dede.Range("K2").Value = "2016-09-01 03:14:15"
$dede.Range("K3").Value = "2016-09-01 03:14:15"

$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
